$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.208.76"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").Value = "3.179.93"
$ws.Range("E3").Value = "  -1.00%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'610.37"
$ws.Range("E5").Value = "  +1.32%  "

$ws.Range("D6").Value = "'155.09"
$ws.Range("E6").Value = "  +0.77%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "3.178.12"
$ws.Range("E8").Value = "  -1.01%  "

$ws.Range("E9").Value = "  +2.45%  "

$ws.Range("E10").Value = "  -1.12%  "

$ws.Range("D11").Value = "'5.67"
$ws.Range("E11").Value = "  -7.47%  "

$ws.Range("D12").Value = "'0.517"
$ws.Range("E12").Value = "  +1.09%  "

$ws.Range("E13").Value = "  -1.23%  "

$ws.Range("D14").Value = "'38.41"
$ws.Range("E14").Value = "  -2.80%  "

$ws.Range("D15").Value = "3.700.88"
$ws.Range("E15").Value = "  -1.03%  "

$ws.Range("D16").Value = "66.243.19"
$ws.Range("E16").Value = "  +0.22%  "

$ws.Range("D17").Value = "'7.42"
$ws.Range("E17").Value = "  -1.14%  "

$ws.Range("D18").Value = "3.181.53"
$ws.Range("E18").Value = "  -0.96%  "

$ws.Range("E19").Value = "  +0.97%  "

$ws.Range("D20").Value = "'510.96"
$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").Value = "'15.41"
$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("D22").Value = "'0.731"
$ws.Range("E22").Value = "  -1.23%  "

$ws.Range("D23").Value = "'8.02"
$ws.Range("E23").Value = "  -0.91%  "

$ws.Range("D24").Value = "'14.82"
$ws.Range("E24").Value = "  -4.80%  "

$ws.Range("D25").Value = "'84.63"
$ws.Range("E25").Value = "  -0.59%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D28").Value = "'9.13"
$ws.Range("E28").Value = "  -2.65%  "

$ws.Range("D29").Value = "'2.38"
$ws.Range("E29").Value = "  +4.17%  "

$ws.Range("B30").Value = "Stacks"
$ws.Range("C30").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D30").Value = "'3.02"
$ws.Range("E30").Value = "  +5.17%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'7.16"
$ws.Range("E31").Value = "  +4.55%  "

$ws.Range("D32").Value = "'28.01"
$ws.Range("E32").Value = "  -0.69%  "

$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("E34").Value = "  -1.81%  "

$ws.Range("D35").Value = "'6.52"
$ws.Range("E35").Value = "  -1.35%  "

$ws.Range("D36").Value = "'505.60"
$ws.Range("E36").Value = "  +4.23%  "

$ws.Range("D37").Value = "'55.11"
$ws.Range("E37").Value = "  -0.25%  "

$ws.Range("D38").Value = "'0.0881"
$ws.Range("E38").Value = "  -2.76%  "

$ws.Range("D39").Value = "'0.0421"
$ws.Range("E39").Value = "  +0.12%  "

$ws.Range("D40").Value = "'0.127"
$ws.Range("E40").Value = "  +5.70%  "

$ws.Range("D41").Value = "'8.80"
$ws.Range("E41").Value = "  -1.47%  "

$ws.Range("D42").Value = "0.0₃0684"
$ws.Range("E42").Value = "  +6.86%  "

$ws.Range("E43").Value = "  -4.33%  "

$ws.Range("E44").Value = "  -2.63%  "

$ws.Range("D45").Value = "'2.43"
$ws.Range("E45").Value = "  -1.13%  "

$ws.Range("D46").Value = "2.828.14"
$ws.Range("E46").Value = "  -4.23%  "

$ws.Range("D47").Value = "'28.09"
$ws.Range("E47").Value = "  -2.19%  "

$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = "  -0.11%  "

$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").Value = "'2.37"
$ws.Range("E49").Value = "  +2.00%  "

$ws.Range("E50").Value = "  +0.39%  "

$ws.Range("E51").Value = "  +6.97%  "
